function Set-TextValue {
    param($range, $value)
    # Force the cell to be written as literal text (matches the workbook's
    # convention of storing even numeric-looking values, e.g. "638.65", as
    # shared strings rather than numbers). Stamping the cell as text first
    # stops Excel's auto-number-detection; resetting the style back to
    # Normal afterwards drops the temporary number-format style so the
    # saved cell carries no explicit style index (matching the source).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Shares": widen column A and append a 3-row gains/losses summary
# ---------------------------------------------------------------------
$wsShares = $wb.Worksheets.Item("Shares")
$wsShares.Columns.Item(1).ColumnWidth = 16.75

Set-TextValue $wsShares.Range("A12") "Gains (incl. losses)"
Set-TextValue $wsShares.Range("I12") "638.65"
Set-TextValue $wsShares.Range("A13") "Gains (excl. losses)"
Set-TextValue $wsShares.Range("I13") "932.75"
Set-TextValue $wsShares.Range("A14") "Losses"
Set-TextValue $wsShares.Range("I14") "-294.10"

# ---------------------------------------------------------------------
# Sheet "Foreign Currencies": widen column A and append the same summary
# ---------------------------------------------------------------------
$wsFx = $wb.Worksheets.Item("Foreign Currencies")
$wsFx.Columns.Item(1).ColumnWidth = 16.75

Set-TextValue $wsFx.Range("A13") "Gains (incl. losses)"
Set-TextValue $wsFx.Range("G13") "15.89"
Set-TextValue $wsFx.Range("A14") "Gains (excl. losses)"
Set-TextValue $wsFx.Range("G14") "77.94"
Set-TextValue $wsFx.Range("A15") "Losses"
Set-TextValue $wsFx.Range("G15") "-62.05"

# ---------------------------------------------------------------------
# Sheet "Dividend Payments": append a Total Amount row
# ---------------------------------------------------------------------
$wsDiv = $wb.Worksheets.Item("Dividend Payments")
Set-TextValue $wsDiv.Range("A5") "Total Amount"
Set-TextValue $wsDiv.Range("E5") "190.67"

# ---------------------------------------------------------------------
# Sheet "Fees": widen column A and append a Total Amount row
# ---------------------------------------------------------------------
$wsFees = $wb.Worksheets.Item("Fees")
$wsFees.Columns.Item(1).ColumnWidth = 12.1

Set-TextValue $wsFees.Range("A12") "Total Amount"
Set-TextValue $wsFees.Range("E12") "29.54"

# ---------------------------------------------------------------------
# Sheet "Tax Withholding": append a Total Amount row
# ---------------------------------------------------------------------
$wsTax = $wb.Worksheets.Item("Tax Withholding")
Set-TextValue $wsTax.Range("A5") "Total Amount"
Set-TextValue $wsTax.Range("E5") "28.60"
